$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new values for column A (Id)
$ws.Range("A2").Value = 109672391
$ws.Range("A3").Value = 109672394
$ws.Range("A4").Value = 109672392

# Set new values for column Q (Ost)
$ws.Range("Q2").Value = 407439.228131063
$ws.Range("Q3").Value = 407928.917497518
$ws.Range("Q4").Value = 407663.9694178636

# Set new values for column R (Nord)
$ws.Range("R2").Value = 7014589.088308994
$ws.Range("R3").Value = 7014334.335588082
$ws.Range("R4").Value = 7014508.711146449

# Set new values for column AC (Publik kommentar)
$ws.Range("AC2").Value = "ringhack"
$ws.Range("AC3").Value = "ringhack gamla"
$ws.Range("AC4").Value = "ringhack gamla"
